$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Design.SlideMaster
$master.Shapes.Item(1).TextFrame.TextRange.Text = $master.Shapes.Item(1).TextFrame.TextRange.Text
Write-Host "done"
